# Auto-generated Excel COM-interop edit script
# Updates the cryptos price/volume table (rows 2-51) to match the
# latest scrape, including two coin re-ranks (rows 17/18, 30/31, 33/34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.762.64'
$ws.Range('E2').Value = '  -2.91%  '

$ws.Range('D3').Value = '3.171.66'
$ws.Range('E3').Value = '  -1.93%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.12'
$ws.Range('E5').Value = '  -1.11%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.71'

$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').Value = '3.165.14'
$ws.Range('E8').Value = '  -2.13%  '

$ws.Range('E9').Value = '  -3.08%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.154'
$ws.Range('E10').Value = '  -4.71%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.63'
$ws.Range('E11').Value = '  -1.72%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.478'
$ws.Range('E12').Value = '  -5.37%  '

$ws.Range('E13').Value = '  -3.91%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.09'
$ws.Range('E14').Value = '  -4.78%  '

$ws.Range('D15').Value = '3.680.20'
$ws.Range('E15').Value = '  -2.24%  '

$ws.Range('D16').Value = '64.847.14'
$ws.Range('E16').Value = '  -2.84%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.179.95'
$ws.Range('E17').Value = '  -1.68%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.114'
$ws.Range('E18').Value = '  +0.56%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.08'
$ws.Range('E19').Value = '  -4.18%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '484.79'
$ws.Range('E20').Value = '  -5.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.85'
$ws.Range('E21').Value = '  -2.51%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.717'
$ws.Range('E22').Value = '  -2.53%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.82'
$ws.Range('E23').Value = '  -2.53%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.96'
$ws.Range('E24').Value = '  -4.93%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.53'
$ws.Range('E25').Value = '  +0.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.14%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.93'
$ws.Range('E27').Value = '  -2.60%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.76'
$ws.Range('E28').Value = '  -4.27%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').Value = '  -4.33%  '

$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.121'
$ws.Range('E30').Value = '  -6.35%  '

$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.12'
$ws.Range('E31').Value = '  +1.65%  '

$ws.Range('E32').Value = '  -7.54%  '

$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.09%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.97'
$ws.Range('E34').Value = '  -4.46%  '

$ws.Range('E35').Value = '  -5.32%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.15'
$ws.Range('E36').Value = '  -5.56%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.69'
$ws.Range('E37').Value = '  -1.67%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.28'
$ws.Range('E38').Value = '  +6.34%  '

$ws.Range('D39').Value = '0.0₃0751'
$ws.Range('E39').Value = '  -2.48%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '458.56'
$ws.Range('E40').Value = '  -9.25%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.127'
$ws.Range('E41').Value = '  -1.89%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0405'
$ws.Range('E42').Value = '  -4.09%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.56'
$ws.Range('E43').Value = '  -2.28%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.45'
$ws.Range('E44').Value = '  -0.16%  '

$ws.Range('D45').Value = '2.908.11'
$ws.Range('E45').Value = '  -0.21%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.277'
$ws.Range('E46').Value = '  -7.55%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.18'
$ws.Range('E47').Value = '  -3.60%  '

$ws.Range('E48').Value = '  -0.02%  '

$ws.Range('E49').Value = '  -2.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.116'
$ws.Range('E50').Value = '  -0.03%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.23'
$ws.Range('E51').Value = '  -1.95%  '
